# Scheduled runner update: refresh market-board sourced profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across the
# leve-profit sheets, per the latest pricing pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 34000
$ws.Range("I86").Value = 60000
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 60000
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -58877
$ws.Range("N86").Value = -10246

$ws.Range("H89").Value = 34000
$ws.Range("I89").Value = 60000
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 300000
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -294384
$ws.Range("N89").Value = -51232

$ws.Range("H101").Value = 1141.875
$ws.Range("J101").Value = 2054.6667
$ws.Range("L101").Value = 6164.000100000001
$ws.Range("N101").Value = -9408.000100000001

$ws.Range("H111").Value = 1158.1666
$ws.Range("I111").Value = 1158.1666
$ws.Range("K111").Value = 3474.4998
$ws.Range("M111").Value = -407.4998000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 227.33333
$ws.Range("I2").Value = 202.8
$ws.Range("K2").Value = 202.8
$ws.Range("M2").Value = -89.80000000000001

$ws.Range("H22").Value = 4333
$ws.Range("J22").Value = 4333
$ws.Range("L22").Value = 4333
$ws.Range("N22").Value = -4931

$ws.Range("H44").Value = 34997.145
$ws.Range("J44").Value = 34997.145
$ws.Range("L44").Value = 34997.145
$ws.Range("N44").Value = -35973.145

$ws.Range("H55").Value = 24998.75
$ws.Range("J55").Value = 24998.75
$ws.Range("L55").Value = 24998.75
$ws.Range("N55").Value = -25628.75

$ws.Range("H61").Value = 1201
$ws.Range("I61").Value = 913.75
$ws.Range("K61").Value = 913.75
$ws.Range("M61").Value = -701.75

$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

$ws.Range("H98").Value = 59996.5
$ws.Range("J98").Value = 59996.5
$ws.Range("L98").Value = 59996.5
$ws.Range("N98").Value = -65986.5

$ws.Range("H116").Value = 227.33333
$ws.Range("I116").Value = 202.8
$ws.Range("K116").Value = 202.8
$ws.Range("M116").Value = 2091.2

$ws.Range("H122").Value = 2997
$ws.Range("I122").Value = 2997
$ws.Range("K122").Value = 8991
$ws.Range("M122").Value = -6541

$ws.Range("H136").Value = 1201
$ws.Range("I136").Value = 913.75
$ws.Range("K136").Value = 2741.25
$ws.Range("M136").Value = -191.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 227.33333
$ws.Range("I3").Value = 202.8
$ws.Range("K3").Value = 202.8
$ws.Range("M3").Value = -88.80000000000001

$ws.Range("H86").Value = 3644.2222
$ws.Range("I86").Value = 3399.8572
$ws.Range("J86").Value = 4499.5
$ws.Range("K86").Value = 3399.8572
$ws.Range("L86").Value = 4499.5
$ws.Range("M86").Value = -2276.8572
$ws.Range("N86").Value = -6745.5

$ws.Range("H89").Value = 3644.2222
$ws.Range("I89").Value = 3399.8572
$ws.Range("J89").Value = 4499.5
$ws.Range("K89").Value = 16999.286
$ws.Range("L89").Value = 22497.5
$ws.Range("M89").Value = -11383.286
$ws.Range("N89").Value = -33729.5

$ws.Range("H100").Value = 31141
$ws.Range("J100").Value = 31141
$ws.Range("L100").Value = 31141
$ws.Range("N100").Value = -33305

$ws.Range("H105").Value = 2081.3333
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 2397.6
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 2397.6
$ws.Range("M105").Value = 1247
$ws.Range("N105").Value = -5891.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 515.8
$ws.Range("J22").Value = 726.6667
$ws.Range("L22").Value = 726.6667
$ws.Range("N22").Value = -1426.6667

$ws.Range("H53").Value = 37499
$ws.Range("J53").Value = 37499
$ws.Range("L53").Value = 37499
$ws.Range("N53").Value = -38713

$ws.Range("H59").Value = 32038.727
$ws.Range("J59").Value = 34713.668
$ws.Range("L59").Value = 34713.668
$ws.Range("N59").Value = -37003.668

$ws.Range("H60").Value = 18782.285
$ws.Range("J60").Value = 24998
$ws.Range("L60").Value = 24998
$ws.Range("N60").Value = -26020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1465.5
$ws.Range("I136").Value = 1465.5
$ws.Range("K136").Value = 4396.5
$ws.Range("M136").Value = 703.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10346

$ws.Range("H24").Value = 17519500
$ws.Range("I24").Value = 35000000
$ws.Range("J24").Value = 39000
$ws.Range("K24").Value = 35000000
$ws.Range("L24").Value = 39000
$ws.Range("M24").Value = -34999827
$ws.Range("N24").Value = -39346

$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10210

$ws.Range("H34").Value = 55000
$ws.Range("J34").Value = 55000
$ws.Range("L34").Value = 55000
$ws.Range("N34").Value = -55536

$ws.Range("H40").Value = 9987.5
$ws.Range("J40").Value = 9987.5
$ws.Range("L40").Value = 9987.5
$ws.Range("N40").Value = -10289.5

$ws.Range("H76").Value = 55000
$ws.Range("J76").Value = 55000
$ws.Range("L76").Value = 55000
$ws.Range("N76").Value = -55630

$ws.Range("H79").Value = 55000
$ws.Range("J79").Value = 55000
$ws.Range("L79").Value = 55000
$ws.Range("N79").Value = -57184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2980.8
$ws.Range("J13").Value = 4682.3335
$ws.Range("L13").Value = 4682.3335
$ws.Range("N13").Value = -4962.3335

$ws.Range("H46").Value = 3998
$ws.Range("J46").Value = 3998
$ws.Range("L46").Value = 3998
$ws.Range("N46").Value = -4374

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 22315.334
$ws.Range("I74").Value = 20500
$ws.Range("J74").Value = 23223
$ws.Range("K74").Value = 20500
$ws.Range("L74").Value = 23223
$ws.Range("M74").Value = -19564
$ws.Range("N74").Value = -25095

$ws.Range("H77").Value = 22315.334
$ws.Range("I77").Value = 20500
$ws.Range("J77").Value = 23223
$ws.Range("K77").Value = 61500
$ws.Range("L77").Value = 69669
$ws.Range("M77").Value = -56820
$ws.Range("N77").Value = -79029

$ws.Range("H96").Value = 1750
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -4746

$ws.Range("H126").Value = 4999
$ws.Range("I126").Value = 4999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12527
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 766.2857
$ws.Range("I132").Value = 694.8333
$ws.Range("K132").Value = 2084.4999
$ws.Range("M132").Value = 445.5001000000002
